$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 1232.6666
$ws.Range("I76").Value = 1299
$ws.Range("J76").Value = 1100
$ws.Range("K76").Value = 1299
$ws.Range("L76").Value = 1100
$ws.Range("M76").Value = -984
$ws.Range("N76").Value = -1730
$ws.Range("H79").Value = 1232.6666
$ws.Range("I79").Value = 1299
$ws.Range("J79").Value = 1100
$ws.Range("K79").Value = 1299
$ws.Range("L79").Value = 1100
$ws.Range("M79").Value = -207
$ws.Range("N79").Value = -3284
$ws.Range("H106").Value = 1751.25
$ws.Range("I106").Value = 1751.25
$ws.Range("K106").Value = 1751.25
$ws.Range("M106").Value = -1120.25
$ws.Range("H107").Value = 298
$ws.Range("I107").Value = 251.78572
$ws.Range("J107").Value = 945
$ws.Range("K107").Value = 251.78572
$ws.Range("L107").Value = 945
$ws.Range("M107").Value = 1668.21428
$ws.Range("N107").Value = -4785
$ws.Range("H125").Value = 41669070
$ws.Range("J125").Value = 3236.25
$ws.Range("L125").Value = 29126.25
$ws.Range("N125").Value = -34046.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5163.533
$ws.Range("I63").Value = 1492.3334
$ws.Range("J63").Value = 7611
$ws.Range("K63").Value = 1492.3334
$ws.Range("L63").Value = 7611
$ws.Range("M63").Value = -806.3334
$ws.Range("N63").Value = -8983
$ws.Range("H66").Value = 5163.533
$ws.Range("I66").Value = 1492.3334
$ws.Range("J66").Value = 7611
$ws.Range("K66").Value = 7461.666999999999
$ws.Range("L66").Value = 38055
$ws.Range("M66").Value = -4029.666999999999
$ws.Range("N66").Value = -44919
$ws.Range("H88").Value = 3503.6875
$ws.Range("I88").Value = 1929
$ws.Range("J88").Value = 4728.4443
$ws.Range("K88").Value = 1929
$ws.Range("L88").Value = 4728.4443
$ws.Range("M88").Value = -1523
$ws.Range("N88").Value = -5540.4443
$ws.Range("H91").Value = 3503.6875
$ws.Range("I91").Value = 1929
$ws.Range("J91").Value = 4728.4443
$ws.Range("K91").Value = 1929
$ws.Range("L91").Value = 4728.4443
$ws.Range("M91").Value = -525
$ws.Range("N91").Value = -7536.4443

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1669.4546
$ws.Range("I20").Value = 1850.5555
$ws.Range("J20").Value = 854.5
$ws.Range("K20").Value = 1850.5555
$ws.Range("L20").Value = 854.5
$ws.Range("M20").Value = -1603.5555
$ws.Range("N20").Value = -1348.5
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = $null
$ws.Range("H94").Value = 694.9091
$ws.Range("I94").Value = 726.6
$ws.Range("K94").Value = 726.6
$ws.Range("M94").Value = -275.6

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19491
$ws.Range("I41").Value = 14321.333
$ws.Range("J41").Value = 35000
$ws.Range("K41").Value = 14321.333
$ws.Range("L41").Value = 35000
$ws.Range("M41").Value = -13893.333
$ws.Range("N41").Value = -35856
$ws.Range("H58").Value = 5451.3125
$ws.Range("I58").Value = 5007.1816
$ws.Range("K58").Value = 5007.1816
$ws.Range("M58").Value = -4804.1816
$ws.Range("H60").Value = 25999.334
$ws.Range("J60").Value = 25999.334
$ws.Range("L60").Value = 25999.334
$ws.Range("N60").Value = -27021.334
$ws.Range("H68").Value = 47649.168
$ws.Range("J68").Value = 47649.168
$ws.Range("L68").Value = 47649.168
$ws.Range("N68").Value = -49147.168
$ws.Range("H71").Value = 47649.168
$ws.Range("J71").Value = 47649.168
$ws.Range("L71").Value = 142947.504
$ws.Range("N71").Value = -150435.504
$ws.Range("H99").Value = 5819.6665
$ws.Range("I99").Value = 5057
$ws.Range("J99").Value = 7802.6
$ws.Range("K99").Value = 5057
$ws.Range("L99").Value = 7802.6
$ws.Range("M99").Value = -3559
$ws.Range("N99").Value = -10798.6
$ws.Range("H122").Value = 1209.2
$ws.Range("I122").Value = 1233.25
$ws.Range("J122").Value = 1113
$ws.Range("K122").Value = 3699.75
$ws.Range("L122").Value = 3339
$ws.Range("M122").Value = -1249.75
$ws.Range("N122").Value = -8239
$ws.Range("H126").Value = 5819.6665
$ws.Range("I126").Value = 5057
$ws.Range("J126").Value = 7802.6
$ws.Range("K126").Value = 15171
$ws.Range("L126").Value = 23407.8
$ws.Range("M126").Value = -12701
$ws.Range("N126").Value = -28347.8
$ws.Range("H132").Value = 9824.4375
$ws.Range("I132").Value = 7529.1
$ws.Range("K132").Value = 22587.3
$ws.Range("M132").Value = -20057.3
$ws.Range("H134").Value = 2301.0417
$ws.Range("I134").Value = 2196.476
$ws.Range("K134").Value = 6589.428
$ws.Range("M134").Value = -4054.428
$ws.Range("H136").Value = 5451.3125
$ws.Range("I136").Value = 5007.1816
$ws.Range("K136").Value = 15021.5448
$ws.Range("M136").Value = -12471.5448

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 44390
$ws.Range("J22").Value = 54150
$ws.Range("L22").Value = 162450
$ws.Range("N22").Value = -162788
$ws.Range("H27").Value = 44390
$ws.Range("J27").Value = 54150
$ws.Range("L27").Value = 162450
$ws.Range("N27").Value = -162654
$ws.Range("H32").Value = 950
$ws.Range("J32").Value = 950
$ws.Range("L32").Value = 2850
$ws.Range("N32").Value = -3416
$ws.Range("H34").Value = 1822.5
$ws.Range("J34").Value = 1822.5
$ws.Range("L34").Value = 5467.5
$ws.Range("N34").Value = -5635.5
$ws.Range("H63").Value = 66
$ws.Range("I63").Value = 66
$ws.Range("K63").Value = 198
$ws.Range("M63").Value = 551
$ws.Range("H64").Value = 2056
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 12000
$ws.Range("N64").Value = -12540
$ws.Range("H66").Value = 66
$ws.Range("I66").Value = 66
$ws.Range("K66").Value = 594
$ws.Range("M66").Value = 3150
$ws.Range("H67").Value = 2056
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 12000
$ws.Range("N67").Value = -13872
$ws.Range("J68").Value = 600
$ws.Range("L68").Value = 1800
$ws.Range("N68").Value = -3422
$ws.Range("J71").Value = 600
$ws.Range("L71").Value = 5400
$ws.Range("N71").Value = -13512
$ws.Range("H76").Value = 825
$ws.Range("I76").Value = 825
$ws.Range("K76").Value = 2475
$ws.Range("M76").Value = -2092
$ws.Range("H79").Value = 825
$ws.Range("I79").Value = 825
$ws.Range("K79").Value = 2475
$ws.Range("M79").Value = -1149
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = $null

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 34000
$ws.Range("J26").Value = 34000
$ws.Range("L26").Value = 34000
$ws.Range("N26").Value = -34560
$ws.Range("H50").Value = 34000
$ws.Range("J50").Value = 34000
$ws.Range("L50").Value = 34000
$ws.Range("N50").Value = -34996
$ws.Range("H122").Value = 3382.8
$ws.Range("I122").Value = 3536.4443
$ws.Range("K122").Value = 10609.3329
$ws.Range("M122").Value = -8159.332900000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 999.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 999.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 999.5
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = -1339.5
$ws.Range("H22").Value = 2166.6667
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -1705
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 2166.6667
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -1893
$ws.Range("N27").Value = -2714
$ws.Range("H46").Value = 2063.5908
$ws.Range("I46").Value = 1628.2858
$ws.Range("J46").Value = 2266.7334
$ws.Range("K46").Value = 1628.2858
$ws.Range("L46").Value = 2266.7334
$ws.Range("M46").Value = -1440.2858
$ws.Range("N46").Value = -2642.7334
$ws.Range("H55").Value = 860.13336
$ws.Range("H82").Value = 2778.5715
$ws.Range("J82").Value = 3263.6365
$ws.Range("L82").Value = 3263.6365
$ws.Range("N82").Value = -3985.6365
$ws.Range("H85").Value = 2778.5715
$ws.Range("J85").Value = 3263.6365
$ws.Range("L85").Value = 3263.6365
$ws.Range("N85").Value = -5759.636500000001
$ws.Range("H100").Value = 3999.3333
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3999.3333
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 3999.3333
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = -5081.3333

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("H126").Value = 1424.2222
$ws.Range("I126").Value = 1424.2222
$ws.Range("K126").Value = 4272.6666
$ws.Range("M126").Value = -1802.6666
$ws.Range("H136").Value = 1393.975
$ws.Range("I136").Value = 1007.91174
$ws.Range("K136").Value = 3023.73522
$ws.Range("M136").Value = -473.73522
